# Append the August/early-September 2021 "Diaria" rows (142-165) that the
# source MV export added to Sheet1, mirroring the existing A:D table
# (Serie / 1 en 1 / 3 en 2 / 5 en 5) built from rows 1-141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r = row, date = column A text, c3/c4 = columns C ("3 en 2") and D ("5 en 5");
# $null means that cell stays empty, exactly like the sparse columns above.
$rows = @(
    @{ r = 142; date = "02-08-2021"; c3 = 2.83; c4 = 2.67 }
    @{ r = 143; date = "03-08-2021"; c3 = 2.78; c4 = 2.75 }
    @{ r = 144; date = "04-08-2021"; c3 = 2.8;  c4 = 2.74 }
    @{ r = 145; date = "05-08-2021"; c3 = 2.78; c4 = 2.91 }
    @{ r = 146; date = "06-08-2021"; c3 = 2.88; c4 = 3.01 }
    @{ r = 147; date = "09-08-2021"; c3 = $null; c4 = 2.97 }
    @{ r = 148; date = "10-08-2021"; c3 = 2.86; c4 = 2.93 }
    @{ r = 149; date = "11-08-2021"; c3 = 2.9;  c4 = 3.04 }
    @{ r = 150; date = "12-08-2021"; c3 = 2.92; c4 = 3 }
    @{ r = 151; date = "13-08-2021"; c3 = $null; c4 = 3.03 }
    @{ r = 152; date = "16-08-2021"; c3 = 3.04; c4 = 3.01 }
    @{ r = 153; date = "17-08-2021"; c3 = $null; c4 = 3.08 }
    @{ r = 154; date = "18-08-2021"; c3 = 2.86; c4 = 3.15 }
    @{ r = 155; date = "19-08-2021"; c3 = 2.83; c4 = $null }
    @{ r = 156; date = "20-08-2021"; c3 = 2.88; c4 = 3.08 }
    @{ r = 157; date = "23-08-2021"; c3 = 2.84; c4 = 3.07 }
    @{ r = 158; date = "24-08-2021"; c3 = $null; c4 = 3.01 }
    @{ r = 159; date = "25-08-2021"; c3 = $null; c4 = 3 }
    @{ r = 160; date = "26-08-2021"; c3 = $null; c4 = 2.98 }
    @{ r = 161; date = "27-08-2021"; c3 = $null; c4 = 3.08 }
    @{ r = 162; date = "31-08-2021"; c3 = 2.91; c4 = 3.11 }
    @{ r = 163; date = "01-09-2021"; c3 = 3.05; c4 = 2.84 }
    @{ r = 164; date = "02-09-2021"; c3 = 2.76; c4 = 3.04 }
    @{ r = 165; date = "03-09-2021"; c3 = 2.81; c4 = 2.97 }
)

foreach ($row in $rows) {
    # Column A dates ("dd-mm-yyyy") must stay plain text, same as the rows
    # above them. Typing a dash-separated date literal makes Excel infer a
    # date serial, so mark the cell Text just for the write, then restore
    # the default "Normal" style so the saved cell carries no special
    # formatting (matching its neighbours exactly).
    $cellA = $ws.Cells.Item($row.r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.date
    $cellA.Style = "Normal"

    if ($null -ne $row.c3) {
        $ws.Cells.Item($row.r, 3).Value = $row.c3
    }
    if ($null -ne $row.c4) {
        $ws.Cells.Item($row.r, 4).Value = $row.c4
    }
}
